$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (row 1) values to lowercase equivalents
$ws.Range("A1").Value = "key"
$ws.Range("B1").Value = "fr"
$ws.Range("C1").Value = "en"

# Row 2 content stays the same (TEST, ceci est un test, that a test)
$ws.Range("A2").Value = "TEST"
$ws.Range("B2").Value = "ceci est un test"
$ws.Range("C2").Value = "that a test"

# Update the active selection to A2 (was A3)
$ws.Range("A2").Select()
